$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Title shape: consolidate the run-per-word text into a single run.
# (Set to a transient placeholder first so the writer doesn't treat the
# assignment as a no-op when the concatenated text already matches.)
$titleRange = $s.Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = " "
$titleRange.Text = "A Table, with a caption"

# Caption textbox: consolidate the run-per-word text into a single run.
$captionRange = $s.Shapes.Item(3).TextFrame.TextRange
$captionRange.Text = " "
$captionRange.Text = "Demonstration of simple table syntax, with alignment"
